$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.189.15'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.04'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.58'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5233'
$ws.Range("E7").Value = '  +2.93%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3774'
$ws.Range("E8").Value = '  +2.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07251'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.21'
$ws.Range("E10").Value = '  +2.84%  '
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07685'
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.904.25'
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '95.31'
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.272'
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008591'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.256.20'
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.058'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.153.45'
$ws.Range("E22").Value = '  +1.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.61'
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.426'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.307'
$ws.Range("E25").Value = '  +10.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.61'
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.14'
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.75'
$ws.Range("E29").Value = '  +1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.964'
$ws.Range("E30").Value = '  +5.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.796'
$ws.Range("E31").Value = '  +2.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09211'
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8077'
$ws.Range("E33").Value = '  +7.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05047'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  +7.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.000'
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.314'
$ws.Range("E37").Value = '  +3.13%  '
$ws.Range("E38").Value = '  +2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5682'
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01983'
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.075'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.977'
$ws.Range("E42").Value = '  +5.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.05'
$ws.Range("E43").Value = '  +3.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.608'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1511'
$ws.Range("E45").Value = '  +1.91%  '
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.16'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.620'
$ws.Range("E49").Value = '  +4.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.50'
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.61'
$ws.Range("E51").Value = '  +0.78%  '
